# Update cryptos list values (prices and 1h volume %) per latest scrape,
# matching the commit 'Updated cryptos list ... with GitHub Actions'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.870.11"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = "'1.895.04"
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'0.7850"
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").Value = "'244.19"
$ws.Range("E6").Value = '  +0.98%  '
$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'0.3152"
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").Value = "'25.46"
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("D10").Value = "'0.07337"
$ws.Range("E10").Value = '  +4.26%  '
$ws.Range("D11").Value = "'0.08130"
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("D12").Value = "'0.7696"
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = "'5.490"
$ws.Range("E13").Value = '  +3.19%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = "'1.792.48"
$ws.Range("E14").Value = '  -5.32%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = "'93.47"
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").Value = "'6.202"
$ws.Range("E16").Value = '  +4.39%  '
$ws.Range("D17").Value = "'29.827.76"
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = "'13.98"
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = "'246.25"
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").Value = "'0.000007869"
$ws.Range("E20").Value = '  +2.25%  '
$ws.Range("D21").Value = "'0.9984"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = "'8.146"
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("D23").Value = "'2.111.26"
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("D24").Value = "'0.9993"
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").Value = "'0.1589"
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("D26").Value = "'9.469"
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("D27").Value = "'162.31"
$ws.Range("E27").Value = '  -1.05%  '
$ws.Range("D28").Value = "'18.81"
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("D29").Value = "'2.035"
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").Value = "'1.468"
$ws.Range("E30").Value = '  +6.41%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").Value = "'4.490"
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = "'0.05607"
$ws.Range("E33").Value = '  -2.15%  '
$ws.Range("D34").Value = "'4.087"
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").Value = "'1.254"
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").Value = "'0.7564"
$ws.Range("E36").Value = '  +2.41%  '
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").Value = "'2.646"
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").Value = "'0.01937"
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = "'2.786"
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").Value = "'1.141.92"
$ws.Range("E41").Value = '  +11.26%  '
$ws.Range("D42").Value = "'0.4467"
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("D43").Value = "'73.73"
$ws.Range("E43").Value = '  +1.80%  '
$ws.Range("D44").Value = "'5.969"
$ws.Range("E44").Value = '  +2.20%  '
$ws.Range("D45").Value = "'0.8568"
$ws.Range("E45").Value = '  +1.91%  '
$ws.Range("D46").Value = "'0.9991"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = "'1.901"
$ws.Range("E47").Value = '  +2.74%  '
$ws.Range("D48").Value = "'3.126"
$ws.Range("E48").Value = '  +6.85%  '
$ws.Range("D49").Value = "'102.08"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").Value = "'9.816"
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").Value = "'7.531"
$ws.Range("E51").Value = '  +0.99%  '
